# Edit script: insert two new data rows (rows 30 and 31) into the daily-price
# log, shifting the existing rows 30-121 down to 32-123, then populate the
# two newly inserted rows with the new observation (Modesto variety, dated
# 2021-12-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 30. This pushes the existing
# rows 30..121 down to 32..123 (matching the target dimension A1:T123),
# while inheriting the formatting (incl. the date-number-format on column D)
# from the surrounding rows.
$ws.Rows("30:31").Insert()

# Row 30: Modesto / Especial, 2021-12-27
$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "Femacal de La Calera"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = "2021-12-27"
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103003
$ws.Range("J30").Value = "Damasco"
$ws.Range("K30").Value = "Modesto"
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 56
$ws.Range("N30").Value = 14000
$ws.Range("O30").Value = 14000
$ws.Range("P30").Value = 14000
$ws.Range("Q30").Value = "$/bandeja 10 kilos"
$ws.Range("R30").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S30").Value = 1400
$ws.Range("T30").Value = 10

# Row 31: Modesto / Primera, 2021-12-27
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = "2021-12-27"
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103003
$ws.Range("J31").Value = "Damasco"
$ws.Range("K31").Value = "Modesto"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 12000
$ws.Range("O31").Value = 12000
$ws.Range("P31").Value = 12000
$ws.Range("Q31").Value = "$/bandeja 10 kilos"
$ws.Range("R31").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S31").Value = 1200
$ws.Range("T31").Value = 10
